$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Double the clockFaceRadius basis (B3) and the pixel basis (N3) so the
#    whole sheet is computed against a 2x bigger radius/pixel scale.
# ---------------------------------------------------------------------------
$ws.Range("B3").Formula = '=2*9.55'
$ws.Range("N3").Value = 260

# ---------------------------------------------------------------------------
# 2. Snapshot the old (pre-change) "Pixels on fr955" O:S values into new
#    columns U:Y so the original numbers stay visible for comparison.
# ---------------------------------------------------------------------------
$ws.Range("U4").Value = 30
$ws.Range("V4").Value = 9
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 92

$ws.Range("U5").Value = 9
$ws.Range("V5").Value = 3
$ws.Range("W5").Value = 3
$ws.Range("X5").Value = 113

$ws.Range("U6").Value = 110
$ws.Range("V6").Value = 16
$ws.Range("W6").Value = 13
$ws.Range("X6").Value = -30

$ws.Range("U7").Value = 146
$ws.Range("V7").Value = 14
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = -31

$ws.Range("U8").Value = 122
$ws.Range("V8").Value = 4
$ws.Range("W8").Value = 4
$ws.Range("X8").Value = -42
$ws.Range("Y8").Value = 13

# ---------------------------------------------------------------------------
# 3. New labels. Set "Pixels based on Original design" (B19) before
#    "Original design" (C10) so the shared-string table grows in the same
#    order as the target workbook (index 16 then 17).
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "Pixels based on Original design"
$ws.Range("H19").Value = "Difference"
$ws.Range("B28").Value = "Pixels based on Original design"
$ws.Range("C10").Value = "Original design"

# ---------------------------------------------------------------------------
# 4. "Original design" measurement table (rows 10-17), mirroring the layout
#    of the "Measurements in centimeters" table up top.
# ---------------------------------------------------------------------------
$ws.Range("C11").Value = "height"
$ws.Range("D11").Value = "width1"
$ws.Range("E11").Value = "width2"
$ws.Range("F11").Value = "radius"
$ws.Range("G11").Value = "circle"

$ws.Range("A12").Value = "clockFaceRadius"
$ws.Range("B12").Value = 100

$ws.Range("A13").Value = "bigTickMark"
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 3.5
$ws.Range("E13").Value = 3.5
$ws.Range("F13").Value = 36.5

$ws.Range("A14").Value = "smallTickMark"
$ws.Range("C14").Value = 3.5
$ws.Range("D14").Value = 1.4
$ws.Range("E14").Value = 1.4
$ws.Range("F14").Value = 45

$ws.Range("A15").Value = "hourHand"
$ws.Range("C15").Value = 44
$ws.Range("D15").Value = 6.3
$ws.Range("E15").Value = 5.0999999999999996
$ws.Range("F15").Value = -12

$ws.Range("A16").Value = "minuteHand"
$ws.Range("C16").Value = 57.8
$ws.Range("D16").Value = 5.2
$ws.Range("E16").Value = 3.7
$ws.Range("F16").Value = -12

$ws.Range("A17").Value = "secondHand"
$ws.Range("C17").Value = 47.9
$ws.Range("D17").Value = 1.4
$ws.Range("E17").Value = 1.4
$ws.Range("F17").Value = -16.5
$ws.Range("G17").Value = 5.0999999999999996

# Cells that are highlighted red in the "Original design" table (style 2).
$ws.Range("D15").Font.Color = 255
$ws.Range("E15").Font.Color = 255
$ws.Range("C16").Font.Color = 255
$ws.Range("E16").Font.Color = 255
$ws.Range("C17").Font.Color = 255
$ws.Range("G17").Font.Color = 255

# ---------------------------------------------------------------------------
# 5. Headers for the "Pixels based on Original design" section (row 20).
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = "height"
$ws.Range("D20").Value = "width1"
$ws.Range("E20").Value = "width2"
$ws.Range("F20").Value = "radius"
$ws.Range("G20").Value = "circle"
$ws.Range("I20").Value = "height"
$ws.Range("J20").Value = "width1"
$ws.Range("K20").Value = "width2"
$ws.Range("L20").Value = "radius"
$ws.Range("M20").Value = "circle"

# ---------------------------------------------------------------------------
# 6. First "Pixels based on Original design" block (rows 21-26), scaled to
#    B21 = 260 pixels, with the 0-decimal ROUND to match the pixel table,
#    and the I:L columns comparing those pixels with the earlier O:R table.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "clockFaceRadius"
$ws.Range("B21").Value = 260

$ws.Range("A22").Value = "bigTickMark"
$ws.Range("C22").Formula = '=ROUND($B$21/$B$12*C13,0)'
$ws.Range("D22:F22").Formula = '=ROUND($B$21/$B$12*D13,0)'
$ws.Range("I22").Formula = '=C22-O13'
$ws.Range("J22:L22").Formula = '=D22-P13'

$ws.Range("A23").Value = "smallTickMark"
$ws.Range("C23:F23").Formula = '=ROUND($B$21/$B$12*C14,0)'
$ws.Range("I23:I26").Formula = '=C23-O14'
$ws.Range("J23:J26").Formula = '=D23-P14'
$ws.Range("K23:K26").Formula = '=E23-Q14'
$ws.Range("L23:M26").Formula = '=F23-R14'

$ws.Range("A24").Value = "hourHand"
$ws.Range("C24:F24").Formula = '=ROUND($B$21/$B$12*C15,0)'

$ws.Range("A25").Value = "minuteHand"
$ws.Range("C25:F25").Formula = '=ROUND($B$21/$B$12*C16,0)'

$ws.Range("A26").Value = "secondHand"
$ws.Range("C26:G26").Formula = '=ROUND($B$21/$B$12*C17,0)'

# ---------------------------------------------------------------------------
# 7. Second "Pixels based on Original design" block (rows 28-35), same idea
#    but rounded to 2 decimals and no Difference columns.
# ---------------------------------------------------------------------------
$ws.Range("C29").Value = "height"
$ws.Range("D29").Value = "width1"
$ws.Range("E29").Value = "width2"
$ws.Range("F29").Value = "radius"
$ws.Range("G29").Value = "circle"

$ws.Range("A30").Value = "clockFaceRadius"
$ws.Range("B30").Value = 260

$ws.Range("A31").Value = "bigTickMark"
$ws.Range("C31").Formula = '=ROUND($B$30/$B$12*C13,2)'
$ws.Range("D31:F31").Formula = '=ROUND($B$30/$B$12*D13,2)'

$ws.Range("A32").Value = "smallTickMark"
$ws.Range("C32:F32").Formula = '=ROUND($B$30/$B$12*C14,2)'

$ws.Range("A33").Value = "hourHand"
$ws.Range("C33:F33").Formula = '=ROUND($B$30/$B$12*C15,2)'

$ws.Range("A34").Value = "minuteHand"
$ws.Range("C34:F34").Formula = '=ROUND($B$30/$B$12*C16,2)'

$ws.Range("A35").Value = "secondHand"
$ws.Range("C35:G35").Formula = '=ROUND($B$30/$B$12*C17,2)'

# ---------------------------------------------------------------------------
# 8. Housekeeping: selection moves to D15 in the refreshed sheet.
# ---------------------------------------------------------------------------
$ws.Range("D15").Select()
